$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row2 = New-Object 'object[,]' 1,10
$row2[0,0] = -0.8011660996346401
$row2[0,1] = 0.5787938608431268
$row2[0,2] = 0.1949980559329599
$row2[0,3] = 0.8798839467072684
$row2[0,4] = 0.1878585623255973
$row2[0,5] = 0.4273069063282261
$row2[0,6] = 0.4640348830873968
$row2[0,7] = 1.290528957474712
$row2[0,8] = 0.6347307473531471
$row2[0,9] = -0.02341949237247376
$ws.Range("B2:K2").Value = $row2

$row3 = New-Object 'object[,]' 1,10
$row3[0,0] = 0.9961641555676001
$row3[0,1] = 1.681050046341908
$row3[0,2] = 0.9890246619602374
$row3[0,3] = 1.228473005962866
$row3[0,4] = 1.265200982722037
$row3[0,5] = 2.091695057109352
$row3[0,6] = 1.435896846987787
$row3[0,7] = 0.7777466072621664
$row3[0,8] = 1.176861555113949
$row3[0,9] = 1.115425012992372
$ws.Range("B3:K3").Value = $row3

$row4 = New-Object 'object[,]' 1,10
$row4[0,0] = -0.007139493607362657
$row4[0,1] = 0.2323088503952662
$row4[0,2] = 0.2690368271544369
$row4[0,3] = 1.095530901541752
$row4[0,4] = 0.4397326914201872
$row4[0,5] = -0.2184175483054337
$row4[0,6] = 0.1806973995463494
$row4[0,7] = 0.119260857424772
$row4[0,8] = -0.3149050776923316
$row4[0,9] = 0.6638963143867453
$ws.Range("B4:K4").Value = $row4

$row5 = New-Object 'object[,]' 1,10
$row5[0,0] = 0.2761763207617995
$row5[0,1] = 1.102670395149115
$row5[0,2] = 0.4468721850275498
$row5[0,3] = -0.211278054698071
$row5[0,4] = 0.187836893153712
$row5[0,5] = 0.1264003510321347
$row5[0,6] = -0.3077655840849689
$row5[0,7] = 0.671035807994108
$row5[0,8] = 0.4439199035461818
$row5[0,9] = -0.1895398258799067
$ws.Range("B5:K5").Value = $row5

$row6 = New-Object 'object[,]' 1,10
$row6[0,0] = 0.1706958642657503
$row6[0,1] = -0.4874543754598706
$row6[0,2] = -0.0883394276080875
$row6[0,3] = -0.1497759697296649
$row6[0,4] = -0.5839419048467684
$row6[0,5] = 0.3948594872323085
$row6[0,6] = 0.1677435827843823
$row6[0,7] = -0.4657161466417062
$row6[0,8] = -0.0383191879861019
$row6[0,9] = -0.2046561426742488
$ws.Range("B6:K6").Value = $row6

$row7 = New-Object 'object[,]' 1,10
$row7[0,0] = -0.2590352918738378
$row7[0,1] = -0.3204718339954152
$row7[0,2] = -0.7546377691125187
$row7[0,3] = 0.2241636229665582
$row7[0,4] = -0.00295228148136796
$row7[0,5] = -0.6364120109074565
$row7[0,6] = -0.2090150522518522
$row7[0,7] = -0.3753520069399991
$row7[0,8] = -0.2526409904991733
$row7[0,9] = -0.205407253944033
$ws.Range("B7:K7").Value = $row7

$row8 = New-Object 'object[,]' 1,10
$row8[0,0] = -0.4956024772386809
$row8[0,1] = 0.483198914840396
$row8[0,2] = 0.2560830103924698
$row8[0,3] = -0.3773767190336187
$row8[0,4] = 0.0500202396219856
$row8[0,5] = -0.1163167150661613
$row8[0,6] = 0.006394301374664513
$row8[0,7] = 0.05362803792980481
$row8[0,8] = -0.3429906935926468
$row8[0,9] = -0.08893865846892673
$ws.Range("B8:K8").Value = $row8

$row9 = New-Object 'object[,]' 1,10
$row9[0,0] = 0.7516854876311507
$row9[0,1] = 0.1182257582050622
$row9[0,2] = 0.5456227168606665
$row9[0,3] = 0.3792857621725196
$row9[0,4] = 0.5019967786133455
$row9[0,5] = 0.5492305151684858
$row9[0,6] = 0.1526117836460341
$row9[0,7] = 0.4066638187697542
$row9[0,8] = 0.6267882086433268
$row9[0,9] = 0.3909997600566587
$ws.Range("B9:K9").Value = $row9

$row10 = New-Object 'object[,]' 1,10
$row10[0,0] = -0.2060627707704842
$row10[0,1] = -0.3723997254586311
$row10[0,2] = -0.2496887090178053
$row10[0,3] = -0.202454972462665
$row10[0,4] = -0.5990737039851166
$row10[0,5] = -0.3450216688613965
$row10[0,6] = -0.1248972789878239
$row10[0,7] = -0.3606857275744921
$row10[0,8] = -0.5396584850452768
$row10[0,9] = -0.5502190153045978
$ws.Range("B10:K10").Value = $row10

$row11 = New-Object 'object[,]' 1,10
$row11[0,0] = -0.04362593824732108
$row11[0,1] = 0.00360779830781921
$row11[0,2] = -0.3930109332146324
$row11[0,3] = -0.1389588980909123
$row11[0,4] = 0.0811654917826603
$row11[0,5] = -0.1546229568040079
$row11[0,6] = -0.3335957142747926
$row11[0,7] = -0.3441562445341136
$row11[0,8] = -0.5518753546922437
$row11[0,9] = -0.334666455217339
$ws.Range("B11:K11").Value = $row11

$row12 = New-Object 'object[,]' 1,10
$row12[0,0] = -0.3493849949673113
$row12[0,1] = -0.09533295984359125
$row12[0,2] = 0.1247914300299814
$row12[0,3] = -0.1109970185566868
$row12[0,4] = -0.2899697760274715
$row12[0,5] = -0.3005303062867926
$row12[0,6] = -0.5082494164449226
$row12[0,7] = -0.2910405169700179
$row12[0,8] = 0.0186456522880436
$row12[0,9] = -0.1356584815439604
$ws.Range("B12:K12").Value = $row12

$row13 = New-Object 'object[,]' 1,10
$row13[0,0] = 0.4741764249972927
$row13[0,1] = 0.2383879764106245
$row13[0,2] = 0.05941521893983981
$row13[0,3] = 0.04885468868051879
$row13[0,4] = -0.1588644214776113
$row13[0,5] = 0.05834447799729348
$row13[0,6] = 0.3680306472553549
$row13[0,7] = 0.2137265134233509
$row13[0,8] = 0.8641529346425341
$row13[0,9] = 0.647987886024325
$ws.Range("B13:K13").Value = $row13

$row14 = New-Object 'object[,]' 1,10
$row14[0,0] = -0.4147612060574529
$row14[0,1] = -0.4253217363167739
$row14[0,2] = -0.633040846474904
$row14[0,3] = -0.4158319469999993
$row14[0,4] = -0.1061457777419378
$row14[0,5] = -0.2604499115739418
$row14[0,6] = 0.3899765096452414
$row14[0,7] = 0.1738114610270322
$row14[0,8] = -0.2332176680079241
$row14[0,9] = 0.2658245491694957
$ws.Range("B14:K14").Value = $row14

$row15 = New-Object 'object[,]' 1,10
$row15[0,0] = -0.2182796404174511
$row15[0,1] = -0.001070740942546333
$row15[0,2] = 0.3086154283155151
$row15[0,3] = 0.1543112944835111
$row15[0,4] = 0.8047377157026943
$row15[0,5] = 0.5885726670844852
$row15[0,6] = 0.1815435380495288
$row15[0,7] = 0.6805857552269486
$row15[0,8] = 0.5217368420714361
$row15[0,9] = 0.5779006630416801
$ws.Range("B15:K15").Value = $row15

$row16 = New-Object 'object[,]' 1,10
$row16[0,0] = 0.5268950687329662
$row16[0,1] = 0.3725909349009622
$row16[0,2] = 1.023017356120145
$row16[0,3] = 0.8068523075019363
$row16[0,4] = 0.3998231784669799
$row16[0,5] = 0.8988653956443997
$row16[0,6] = 0.7400164824888872
$row16[0,7] = 0.7961803034591312
$row16[0,8] = 3.045983143070524
$row16[0,9] = 10.51513347392782
$ws.Range("B16:K16").Value = $row16

$row17 = New-Object 'object[,]' 1,10
$row17[0,0] = 0.1553820354260574
$row17[0,1] = 0.8058084566452406
$row17[0,2] = 0.5896434080270315
$row17[0,3] = 0.1826142789920752
$row17[0,4] = 0.681656496169495
$row17[0,5] = 0.5228075830139824
$row17[0,6] = 0.5789714039842264
$row17[0,7] = 2.828774243595618
$row17[0,8] = 10.29792457445291
$row17[0,9] = -7.873337225684557
$ws.Range("B17:K17").Value = $row17

$row18 = New-Object 'object[,]' 1,10
$row18[0,0] = 0.4961222873871792
$row18[0,1] = 0.2799572387689701
$row18[0,2] = -0.1270718902659863
$row18[0,3] = 0.3719703269114335
$row18[0,4] = 0.2131214137559209
$row18[0,5] = 0.269285234726165
$row18[0,6] = 2.519088074337557
$row18[0,7] = 9.988238405194851
$row18[0,8] = -8.183023394942618
$row18[0,9] = -0.04335184796346914
$ws.Range("B18:K18").Value = $row18

$row19 = New-Object 'object[,]' 1,10
$row19[0,0] = 0.4342613726009741
$row19[0,1] = 0.02723224356601772
$row19[0,2] = 0.5262744607434375
$row19[0,3] = 0.3674255475879249
$row19[0,4] = 0.423589368558169
$row19[0,5] = 2.673392208169561
$row19[0,6] = 10.14254253902685
$row19[0,7] = -8.028719261110615
$row19[0,8] = 0.1109522858685349
$row19[0,9] = 2.258522757667927
$ws.Range("B19:K19").Value = $row19

$row20 = New-Object 'object[,]' 1,10
$row20[0,0] = -0.6231941776531655
$row20[0,1] = -0.1241519604757457
$row20[0,2] = -0.2830008736312583
$row20[0,3] = -0.2268370526610142
$row20[0,4] = 2.022965786950378
$row20[0,5] = 9.492116117807671
$row20[0,6] = -8.679145682329798
$row20[0,7] = -0.5394741353506483
$row20[0,8] = 1.608096336448744
$row20[0,9] = -1.831933650074586
$ws.Range("B20:K20").Value = $row20

$row21 = New-Object 'object[,]' 1,10
$row21[0,0] = 0.09201308814246346
$row21[0,1] = -0.06683582501304909
$row21[0,2] = -0.01067200404280504
$row21[0,3] = 2.239130835568587
$row21[0,4] = 9.708281166425881
$row21[0,5] = -8.462980633711588
$row21[0,6] = -0.3233090867324392
$row21[0,7] = 1.824261385066953
$row21[0,8] = -1.615768601456377
$row21[0,9] = -1.676942779870578
$ws.Range("B21:K21").Value = $row21

$row22 = New-Object 'object[,]' 1,10
$row22[0,0] = 0.3401933040219072
$row22[0,1] = 0.3963571249921513
$row22[0,2] = 2.646159964603544
$row22[0,3] = 10.11531029546084
$row22[0,4] = -8.055951504676631
$row22[0,5] = 0.08372004230251717
$row22[0,6] = 2.231290514101909
$row22[0,7] = -1.20873947242142
$row22[0,8] = -1.269913650835621
$row22[0,9] = 0.7715422863206332
$ws.Range("B22:K22").Value = $row22

$row23 = New-Object 'object[,]' 1,10
$row23[0,0] = -0.1026850921852685
$row23[0,1] = 2.147117747426124
$row23[0,2] = 9.616268078283417
$row23[0,3] = -8.554993721854052
$row23[0,4] = -0.4153221748749026
$row23[0,5] = 1.73224829692449
$row23[0,6] = -1.70778168959884
$row23[0,7] = -1.768955868013041
$row23[0,8] = 0.2725000691432133
$row23[0,9] = -0.3100944678904277
$ws.Range("B23:K23").Value = $row23

$row24 = New-Object 'object[,]' 1,10
$row24[0,0] = 2.305966660581636
$row24[0,1] = 9.77511699143893
$row24[0,2] = -8.396144808698539
$row24[0,3] = -0.2564732617193901
$row24[0,4] = 1.891097210080002
$row24[0,5] = -1.548932776443328
$row24[0,6] = -1.610106954857529
$row24[0,7] = 0.4313489822987259
$row24[0,8] = -0.1512455547349151
$row24[0,9] = -0.05557051392555735
$ws.Range("B24:K24").Value = $row24

